# Auto-generated edit script applying the diff to Jenova_Profits.xlsx
# Updates currentAveragePrice / profit-calculation columns (H-N) on several rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 28
$ws.Range("H28").Value = 68343.53
$ws.Range("I28").Value = 126432.75
$ws.Range("K28").Value = 126432.75
$ws.Range("M28").Value = -125947.75

# Row 43
$ws.Range("H43").Value = 2362.0667
$ws.Range("I43").Value = 3500
$ws.Range("J43").Value = 1366.375
$ws.Range("K43").Value = 3500
$ws.Range("L43").Value = 1366.375
$ws.Range("M43").Value = -3431
$ws.Range("N43").Value = -1504.375

# Row 47
$ws.Range("H47").Value = 4000
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents() | Out-Null

# Row 58
$ws.Range("H58").Value = 7409.647
$ws.Range("J58").Value = 9532.308000000001
$ws.Range("L58").Value = 28596.924
$ws.Range("N58").Value = -28896.924

# Row 69
$ws.Range("H69").Value = 16125
$ws.Range("J69").Value = 14833.333
$ws.Range("L69").Value = 44499.999
$ws.Range("N69").Value = -46247.999

# Row 72
$ws.Range("H72").Value = 16125
$ws.Range("J72").Value = 14833.333
$ws.Range("L72").Value = 133499.997
$ws.Range("N72").Value = -142235.997

# Row 127
$ws.Range("H127").Value = 10335.385
$ws.Range("I127").Value = 12356.3
$ws.Range("K127").Value = 37068.89999999999
$ws.Range("M127").Value = -32108.89999999999

# Row 131
$ws.Range("H131").Value = 2929.4736
$ws.Range("I131").Value = 2576.8572
$ws.Range("J131").Value = 3916.8
$ws.Range("K131").Value = 7730.571599999999
$ws.Range("L131").Value = 11750.4
$ws.Range("M131").Value = -2690.571599999999
$ws.Range("N131").Value = -21830.4

# Row 138
$ws.Range("H138").Value = 4274.356
$ws.Range("J138").Value = 5314.676
$ws.Range("L138").Value = 15944.028
$ws.Range("N138").Value = -26224.028

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 31
$ws.Range("H31").Value = 6999.5
$ws.Range("I31").Value = 6999.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6999.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -6705.5
$ws.Range("N31").ClearContents() | Out-Null

# Row 45
$ws.Range("H45").Value = 6429.231
$ws.Range("I45").Value = 4013.8333
$ws.Range("K45").Value = 4013.8333
$ws.Range("M45").Value = -3636.8333

# Row 74
$ws.Range("H74").Value = 1398.9697
$ws.Range("I74").Value = 939.06665
$ws.Range("K74").Value = 939.06665
$ws.Range("M74").Value = -65.06664999999998

# Row 77
$ws.Range("H77").Value = 1398.9697
$ws.Range("I77").Value = 939.06665
$ws.Range("K77").Value = 4695.33325
$ws.Range("M77").Value = -327.3332499999997

# Row 102
$ws.Range("H102").Value = 2312.35
$ws.Range("I102").Value = 2170.8948
$ws.Range("K102").Value = 2170.8948
$ws.Range("M102").Value = -548.8948

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 94
$ws.Range("H94").Value = 1999.3334
$ws.Range("J94").Value = 1999.5
$ws.Range("L94").Value = 1999.5
$ws.Range("N94").Value = -2901.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 19
$ws.Range("H19").Value = 2427.0715
$ws.Range("I19").Value = 163.16667
$ws.Range("K19").Value = 163.16667
$ws.Range("M19").Value = 6.833329999999989

# Row 20
$ws.Range("H20").Value = 71948
$ws.Range("J20").Value = 71948
$ws.Range("L20").Value = 71948
$ws.Range("N20").Value = -72420

# Row 22
$ws.Range("H22").Value = 466
$ws.Range("J22").Value = 825
$ws.Range("L22").Value = 825
$ws.Range("N22").Value = -1525

# Row 24
$ws.Range("H24").Value = 2427.0715
$ws.Range("I24").Value = 163.16667
$ws.Range("K24").Value = 163.16667
$ws.Range("M24").Value = 6.833329999999989

# Row 30
$ws.Range("H30").Value = 71948
$ws.Range("J30").Value = 71948
$ws.Range("L30").Value = 71948
$ws.Range("N30").Value = -72130

# Row 37
$ws.Range("H37").Value = 7799.4
$ws.Range("J37").Value = 7799.4
$ws.Range("L37").Value = 7799.4
$ws.Range("N37").Value = -8013.4

# Row 56
$ws.Range("H56").Value = 12000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents() | Out-Null

# Row 128
$ws.Range("H128").Value = 71948
$ws.Range("J128").Value = 71948
$ws.Range("L128").Value = 71948
$ws.Range("N128").Value = -81908

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 2
$ws.Range("H2").Value = 89.36364
$ws.Range("I2").Value = 34.22222
$ws.Range("K2").Value = 205.33332
$ws.Range("M2").Value = -92.33332000000001

# Row 14
$ws.Range("H14").Value = 2278.3
$ws.Range("I14").Value = 2278.3
$ws.Range("K14").Value = 6834.900000000001
$ws.Range("M14").Value = -6661.900000000001

# Row 132
$ws.Range("H132").Value = 3052.8125
$ws.Range("I132").Value = 2300
$ws.Range("J132").Value = 3447.1428
$ws.Range("K132").Value = 20700
$ws.Range("L132").Value = 31024.2852
$ws.Range("M132").Value = -18170
$ws.Range("N132").Value = -36084.2852

# Row 141
$ws.Range("H141").Value = 10371.6
$ws.Range("I141").Value = 10530.857
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 31592.571
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -26412.571
$ws.Range("N141").Value = -40360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 7149.645
$ws.Range("I102").Value = 6982.8
$ws.Range("K102").Value = 6982.8
$ws.Range("M102").Value = -5360.8

# Row 107
$ws.Range("H107").Value = 1079.2
$ws.Range("J107").Value = 2198
$ws.Range("L107").Value = 2198
$ws.Range("N107").Value = -6038

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 2107.7778
$ws.Range("I22").Value = 1284.8572
$ws.Range("K22").Value = 1284.8572
$ws.Range("M22").Value = -989.8571999999999

# Row 27
$ws.Range("H27").Value = 2107.7778
$ws.Range("I27").Value = 1284.8572
$ws.Range("K27").Value = 1284.8572
$ws.Range("M27").Value = -1177.8572

# Row 40
$ws.Range("H40").Value = 3004733
$ws.Range("I40").Value = 5458124.5
$ws.Range("J40").Value = 6143.222
$ws.Range("K40").Value = 5458124.5
$ws.Range("L40").Value = 6143.222
$ws.Range("M40").Value = -5457988.5
$ws.Range("N40").Value = -6415.222

# Row 46
$ws.Range("H46").Value = 3467.818
$ws.Range("J46").Value = 4500.9
$ws.Range("L46").Value = 4500.9
$ws.Range("N46").Value = -4876.9

# Row 68
$ws.Range("H68").Value = 9238
$ws.Range("I68").Value = 9085.799999999999
$ws.Range("J68").Value = 9999
$ws.Range("K68").Value = 9085.799999999999
$ws.Range("L68").Value = 9999
$ws.Range("M68").Value = -8336.799999999999
$ws.Range("N68").Value = -11497

# Row 71
$ws.Range("H71").Value = 9238
$ws.Range("I71").Value = 9085.799999999999
$ws.Range("J71").Value = 9999
$ws.Range("K71").Value = 45429
$ws.Range("L71").Value = 49995
$ws.Range("M71").Value = -41685
$ws.Range("N71").Value = -57483

# Row 82
$ws.Range("H82").Value = 4240.4
$ws.Range("I82").Value = 4240.4
$ws.Range("K82").Value = 4240.4
$ws.Range("M82").Value = -3879.4

# Row 85
$ws.Range("H85").Value = 4240.4
$ws.Range("I85").Value = 4240.4
$ws.Range("K85").Value = 4240.4
$ws.Range("M85").Value = -2992.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 62
$ws.Range("H62").Value = 9489.556
$ws.Range("J62").Value = 9841.200000000001
$ws.Range("L62").Value = 9841.200000000001
$ws.Range("N62").Value = -11089.2

# Row 65
$ws.Range("H65").Value = 9489.556
$ws.Range("J65").Value = 9841.200000000001
$ws.Range("L65").Value = 49206
$ws.Range("N65").Value = -55446
